# Complete Unit Test for BookLendingApiController
#
# Row 6 of the plan ("Add more client pages to use APIs (Add new User/Book)")
# moves from "In-Progress" to "Completed":
#   - F6 (Actual End-date) gets an actual completion date (2019-09-19 / serial 43727),
#     formatted/bordered the same way as the other "Actual End-date" cells (F3:F5).
#   - G6 (Status) switches from the one-off green "In-Progress" text/style to the
#     same bold "Completed" style used by the other finished rows (G3:G5).
# The active selection also advances from F6 to F7, matching where the author's
# cursor ended up after filling the row in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- F6: stamp the actual end-date, copying the date formatting already used
#         by the sibling "Actual End-date" cells (F3) so the new cell gets the
#         same number format / borders rather than the old blank style.
$ws.Cells.Item(3, 6).Copy()
$ws.Cells.Item(6, 6).PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(6, 6).Value = 43727          # 2019-09-19

# --- G6: mark the row "Completed", reusing the formatting already used by
#         the other completed rows (G3) instead of the removed green font.
$ws.Cells.Item(3, 7).Copy()
$ws.Cells.Item(6, 7).PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(6, 7).Value = "Completed"

$excel.CutCopyMode = 0

# --- Move the active selection to F7, where the author's cursor landed next.
$ws.Range("F7").Select()
